$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.334.81'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.638.91'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.57%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.22'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.529'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.29'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.256'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0608'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.869.00'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.640.18'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.01'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.557'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.19'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.301.84'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '227.37'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -8.51%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0719'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.36'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.28'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.57%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.25'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '146.84'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.93'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.51'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.68%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.15%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0484'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -5.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.26'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.63%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.07'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.393.68'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.63%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.55'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.876'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -6.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.553'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0165'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.02'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.87%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('B42').Value = 'mCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.48'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.43'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.21'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.782'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.78%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '63.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.17%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.780.02'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.64'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.96%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '86.71'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.02%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₆0105'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.17%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0981'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.72%  '
